{"js": "const NCOLS = 5;\nconst pairs = [\n  [\"71\u00d753=\", \"44\u00d760=\"],\n  [\"74\u00d767=\", \"70\u00d745=\"],\n  [\"84\u00d725=\", \"32\u00d785=\"],\n  [\"79\u00d792=\", \"87\u00d733=\"],\n  [\"68\u00d710=\", \"64\u00d738=\"],\n  [\"26\u00d794=\", \"60\u00d795=\"],\n  [\"25\u00d794=\", \"15\u00d779=\"],\n  [\"19\u00d789=\", \"17\u00d721=\"],\n  [\"44\u00d751=\", \"34\u00d721=\"],\n  [\"38\u00d745=\", \"35\u00d747=\"],\n  [\"42\u00d715=\", \"13\u00d740=\"],\n  [\"64\u00d732=\", \"20\u00d732=\"],\n  [\"16\u00d723=\", \"53\u00d780=\"],\n  [\"13\u00d764=\", \"37\u00d794=\"],\n  [\"100\u00d768=\", \"87\u00d773=\"],\n  [\"17\u00d759=\", \"22\u00d776=\"],\n  [\"59\u00d740=\", \"56\u00d772=\"],\n  [\"66\u00d735=\", \"76\u00d774=\"],\n  [\"45\u00d789=\", \"69\u00d728=\"],\n  [\"91\u00d744=\", \"70\u00d759=\"],\n  [\"64\u00d719=\", \"43\u00d767=\"],\n  [\"71\u00d728=\", \"73\u00d753=\"],\n  [\"55\u00d776=\", \"13\u00d744=\"],\n  [\"65\u00d727=\", \"23\u00d744=\"],\n  [\"21\u00d747=\", \"10\u00d713=\"],\n  [\"63\u00d758=\", \"98\u00d757=\"],\n  [\"99\u00d796=\", \"83\u00d756=\"],\n  [\"36\u00d780=\", \"17\u00d732=\"],\n  [\"82\u00d715=\", \"79\u00d776=\"],\n  [\"37\u00d755=\", \"11\u00d786=\"],\n  [\"68\u00d788=\", \"43\u00d780=\"],\n  [\"48\u00d772=\", \"64\u00d728=\"],\n  [\"56\u00d723=\", \"32\u00d758=\"],\n  [\"49\u00d729=\", \"92\u00d774=\"],\n  [\"76\u00d725=\", \"12\u00d771=\"],\n  [\"12\u00d777=\", \"23\u00d781=\"],\n  [\"34\u00d747=\", \"50\u00d788=\"],\n  [\"54\u00d793=\", \"74\u00d769=\"],\n  [\"32\u00d773=\", \"76\u00d713=\"],\n  [\"84\u00d725=\", \"27\u00d720=\"],\n  [\"59\u00d728=\", \"31\u00d712=\"],\n  [\"70\u00d735=\", \"89\u00d765=\"],\n  [\"50\u00d746=\", \"33\u00d740=\"],\n  [\"96\u00d794=\", \"50\u00d721=\"],\n  [\"80\u00d7100=\", \"71\u00d754=\"],\n  [\"79\u00d747=\", \"78\u00d716=\"],\n  [\"51\u00d749=\", \"44\u00d770=\"],\n  [\"61\u00d719=\", \"83\u00d776=\"],\n  [\"44\u00d796=\", \"63\u00d774=\"],\n  [\"84\u00d755=\", \"97\u00d723=\"],\n  [\"38\u00d738=\", \"47\u00d762=\"],\n  [\"93\u00d736=\", \"99\u00d780=\"],\n  [\"65\u00d714=\", \"57\u00d790=\"],\n  [\"16\u00d797=\", \"10\u00d746=\"],\n  [\"35\u00d719=\", \"90\u00d723=\"],\n  [\"86\u00d774=\", \"16\u00d727=\"],\n  [\"76\u00d745=\", \"52\u00d739=\"],\n  [\"61\u00d739=\", \"87\u00d752=\"],\n  [\"12\u00d776=\", \"15\u00d771=\"],\n  [\"72\u00d783=\", \"37\u00d712=\"],\n  [\"100\u00d714=\", \"90\u00d740=\"],\n  [\"90\u00d753=\", \"42\u00d743=\"],\n  [\"15\u00d778=\", \"30\u00d760=\"],\n  [\"97\u00d782=\", \"84\u00d775=\"],\n  [\"22\u00d787=\", \"22\u00d761=\"],\n  [\"22\u00d745=\", \"16\u00d717=\"],\n  [\"18\u00d754=\", \"83\u00d777=\"],\n  [\"61\u00d717=\", \"78\u00d764=\"],\n  [\"22\u00d747=\", \"24\u00d748=\"],\n  [\"62\u00d714=\", \"100\u00d735=\"],\n  [\"33\u00d780=\", \"14\u00d740=\"],\n  [\"55\u00d729=\", \"84\u00d798=\"],\n  [\"84\u00d711=\", \"45\u00d729=\"],\n  [\"90\u00d762=\", \"85\u00d715=\"],\n  [\"33\u00d779=\", \"72\u00d780=\"],\n  [\"19\u00d715=\", \"53\u00d729=\"],\n  [\"80\u00d722=\", \"92\u00d742=\"],\n  [\"80\u00d734=\", \"19\u00d775=\"],\n  [\"46\u00d748=\", \"36\u00d749=\"],\n  [\"55\u00d783=\", \"54\u00d733=\"],\n  [\"77\u00d747=\", \"12\u00d773=\"],\n  [\"51\u00d783=\", \"14\u00d760=\"],\n  [\"11\u00d768=\", \"99\u00d733=\"],\n  [\"55\u00d759=\", \"58\u00d776=\"],\n  [\"34\u00d758=\", \"92\u00d773=\"],\n  [\"53\u00d765=\", \"98\u00d745=\"],\n  [\"30\u00d798=\", \"66\u00d730=\"],\n  [\"35\u00d788=\", \"22\u00d726=\"],\n  [\"36\u00d737=\", \"70\u00d785=\"],\n  [\"11\u00d732=\", \"97\u00d761=\"],\n  [\"39\u00d780=\", \"94\u00d756=\"],\n  [\"27\u00d738=\", \"62\u00d726=\"],\n  [\"27\u00d775=\", \"56\u00d790=\"],\n  [\"71\u00d798=\", \"77\u00d735=\"],\n  [\"88\u00d793=\", \"54\u00d765=\"],\n  [\"86\u00d745=\", \"15\u00d787=\"],\n  [\"32\u00d726=\", \"86\u00d767=\"],\n  [\"31\u00d767=\", \"40\u00d789=\"],\n  [\"53\u00d770=\", \"96\u00d748=\"],\n  [\"54\u00d789=\", \"88\u00d717=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (let i = 0; i < pairs.length; i++) {\n  const row = Math.floor(i / NCOLS);\n  const col = i % NCOLS;\n  const [oldText, newText] = pairs[i];\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`cell (${row},${col}) expected text \"${oldText}\" not found`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$cols = 5\n$pairs = @(\n  ,@(\"71\u00d753=\", \"44\u00d760=\")\n  ,@(\"74\u00d767=\", \"70\u00d745=\")\n  ,@(\"84\u00d725=\", \"32\u00d785=\")\n  ,@(\"79\u00d792=\", \"87\u00d733=\")\n  ,@(\"68\u00d710=\", \"64\u00d738=\")\n  ,@(\"26\u00d794=\", \"60\u00d795=\")\n  ,@(\"25\u00d794=\", \"15\u00d779=\")\n  ,@(\"19\u00d789=\", \"17\u00d721=\")\n  ,@(\"44\u00d751=\", \"34\u00d721=\")\n  ,@(\"38\u00d745=\", \"35\u00d747=\")\n  ,@(\"42\u00d715=\", \"13\u00d740=\")\n  ,@(\"64\u00d732=\", \"20\u00d732=\")\n  ,@(\"16\u00d723=\", \"53\u00d780=\")\n  ,@(\"13\u00d764=\", \"37\u00d794=\")\n  ,@(\"100\u00d768=\", \"87\u00d773=\")\n  ,@(\"17\u00d759=\", \"22\u00d776=\")\n  ,@(\"59\u00d740=\", \"56\u00d772=\")\n  ,@(\"66\u00d735=\", \"76\u00d774=\")\n  ,@(\"45\u00d789=\", \"69\u00d728=\")\n  ,@(\"91\u00d744=\", \"70\u00d759=\")\n  ,@(\"64\u00d719=\", \"43\u00d767=\")\n  ,@(\"71\u00d728=\", \"73\u00d753=\")\n  ,@(\"55\u00d776=\", \"13\u00d744=\")\n  ,@(\"65\u00d727=\", \"23\u00d744=\")\n  ,@(\"21\u00d747=\", \"10\u00d713=\")\n  ,@(\"63\u00d758=\", \"98\u00d757=\")\n  ,@(\"99\u00d796=\", \"83\u00d756=\")\n  ,@(\"36\u00d780=\", \"17\u00d732=\")\n  ,@(\"82\u00d715=\", \"79\u00d776=\")\n  ,@(\"37\u00d755=\", \"11\u00d786=\")\n  ,@(\"68\u00d788=\", \"43\u00d780=\")\n  ,@(\"48\u00d772=\", \"64\u00d728=\")\n  ,@(\"56\u00d723=\", \"32\u00d758=\")\n  ,@(\"49\u00d729=\", \"92\u00d774=\")\n  ,@(\"76\u00d725=\", \"12\u00d771=\")\n  ,@(\"12\u00d777=\", \"23\u00d781=\")\n  ,@(\"34\u00d747=\", \"50\u00d788=\")\n  ,@(\"54\u00d793=\", \"74\u00d769=\")\n  ,@(\"32\u00d773=\", \"76\u00d713=\")\n  ,@(\"84\u00d725=\", \"27\u00d720=\")\n  ,@(\"59\u00d728=\", \"31\u00d712=\")\n  ,@(\"70\u00d735=\", \"89\u00d765=\")\n  ,@(\"50\u00d746=\", \"33\u00d740=\")\n  ,@(\"96\u00d794=\", \"50\u00d721=\")\n  ,@(\"80\u00d7100=\", \"71\u00d754=\")\n  ,@(\"79\u00d747=\", \"78\u00d716=\")\n  ,@(\"51\u00d749=\", \"44\u00d770=\")\n  ,@(\"61\u00d719=\", \"83\u00d776=\")\n  ,@(\"44\u00d796=\", \"63\u00d774=\")\n  ,@(\"84\u00d755=\", \"97\u00d723=\")\n  ,@(\"38\u00d738=\", \"47\u00d762=\")\n  ,@(\"93\u00d736=\", \"99\u00d780=\")\n  ,@(\"65\u00d714=\", \"57\u00d790=\")\n  ,@(\"16\u00d797=\", \"10\u00d746=\")\n  ,@(\"35\u00d719=\", \"90\u00d723=\")\n  ,@(\"86\u00d774=\", \"16\u00d727=\")\n  ,@(\"76\u00d745=\", \"52\u00d739=\")\n  ,@(\"61\u00d739=\", \"87\u00d752=\")\n  ,@(\"12\u00d776=\", \"15\u00d771=\")\n  ,@(\"72\u00d783=\", \"37\u00d712=\")\n  ,@(\"100\u00d714=\", \"90\u00d740=\")\n  ,@(\"90\u00d753=\", \"42\u00d743=\")\n  ,@(\"15\u00d778=\", \"30\u00d760=\")\n  ,@(\"97\u00d782=\", \"84\u00d775=\")\n  ,@(\"22\u00d787=\", \"22\u00d761=\")\n  ,@(\"22\u00d745=\", \"16\u00d717=\")\n  ,@(\"18\u00d754=\", \"83\u00d777=\")\n  ,@(\"61\u00d717=\", \"78\u00d764=\")\n  ,@(\"22\u00d747=\", \"24\u00d748=\")\n  ,@(\"62\u00d714=\", \"100\u00d735=\")\n  ,@(\"33\u00d780=\", \"14\u00d740=\")\n  ,@(\"55\u00d729=\", \"84\u00d798=\")\n  ,@(\"84\u00d711=\", \"45\u00d729=\")\n  ,@(\"90\u00d762=\", \"85\u00d715=\")\n  ,@(\"33\u00d779=\", \"72\u00d780=\")\n  ,@(\"19\u00d715=\", \"53\u00d729=\")\n  ,@(\"80\u00d722=\", \"92\u00d742=\")\n  ,@(\"80\u00d734=\", \"19\u00d775=\")\n  ,@(\"46\u00d748=\", \"36\u00d749=\")\n  ,@(\"55\u00d783=\", \"54\u00d733=\")\n  ,@(\"77\u00d747=\", \"12\u00d773=\")\n  ,@(\"51\u00d783=\", \"14\u00d760=\")\n  ,@(\"11\u00d768=\", \"99\u00d733=\")\n  ,@(\"55\u00d759=\", \"58\u00d776=\")\n  ,@(\"34\u00d758=\", \"92\u00d773=\")\n  ,@(\"53\u00d765=\", \"98\u00d745=\")\n  ,@(\"30\u00d798=\", \"66\u00d730=\")\n  ,@(\"35\u00d788=\", \"22\u00d726=\")\n  ,@(\"36\u00d737=\", \"70\u00d785=\")\n  ,@(\"11\u00d732=\", \"97\u00d761=\")\n  ,@(\"39\u00d780=\", \"94\u00d756=\")\n  ,@(\"27\u00d738=\", \"62\u00d726=\")\n  ,@(\"27\u00d775=\", \"56\u00d790=\")\n  ,@(\"71\u00d798=\", \"77\u00d735=\")\n  ,@(\"88\u00d793=\", \"54\u00d765=\")\n  ,@(\"86\u00d745=\", \"15\u00d787=\")\n  ,@(\"32\u00d726=\", \"86\u00d767=\")\n  ,@(\"31\u00d767=\", \"40\u00d789=\")\n  ,@(\"53\u00d770=\", \"96\u00d748=\")\n  ,@(\"54\u00d789=\", \"88\u00d717=\")\n)\n\nfor ($i = 0; $i -lt $pairs.Count; $i++) {\n  $row = [int][math]::Floor($i / $cols) + 1\n  $col = ($i % $cols) + 1\n  $oldText = $pairs[$i][0]\n  $newText = $pairs[$i][1]\n  $cell = $tbl.Cell($row, $col)\n  $checkRng = $cell.Range\n  [void]$checkRng.MoveEnd(1, -1)\n  if ($checkRng.Text -ne $oldText) {\n    throw \"cell ($row,$col) expected text [$oldText] but found [$($checkRng.Text)]\"\n  }\n  $cell.Range.Text = $newText\n}\nWrite-Output \"done: updated $($pairs.Count) cells\""}
